$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[49.944756189191175, 50.06246561035622]"
$ws.Range("T2").Value = "[49.93257738792698, 50.015411190681476]"
$ws.Range("L3").Value = "[49.98731326050781, 50.175324133080444]"
$ws.Range("T3").Value = "[49.97629183600236, 50.073932859137244]"
